# heat delivery set allocation, energy flow updates
$wb = $excel.ActiveWorkbook

$wsConsumption  = $wb.Worksheets.Item(1)   # consumptionAssets
$wsProduction   = $wb.Worksheets.Item(2)   # productionAssets
$wsConversion   = $wb.Worksheets.Item(3)   # conversionAssets
$wsStorage      = $wb.Worksheets.Item(4)   # storageAssets

# ---------------------------------------------------------------------------
# conversionAssets (sheet3): rename the HEATDELIVERYSET energyAssetType and
# shrink its electricity capacity; widen the name column.
# ---------------------------------------------------------------------------
$wsConversion.Range("D6").Value = "HEAT_DELIVERY_SET"
$wsConversion.Range("F6").Value = 10

$wsConversion.Columns.Item(2).ColumnWidth = 33.4

# ---------------------------------------------------------------------------
# storageAssets (sheet4): separate calculation / operation of district
# heating netconnections - rework capacity formulas (kW -> MW, /1000) and
# bump several allocation / capacity figures.
#
# Number formats are applied in the same order the original authoring did
# (scientific notation for the big L-column capacities first, then the one
# decimal place format for the re-based F-column capacities) so the
# generated cellXfs/numFmts indices line up: s="1" -> 0.00E+00 (builtin 11),
# s="2" -> custom "0.0" (164).
# ---------------------------------------------------------------------------

$wsStorage.Range("L5:L9").NumberFormat = "0.00E+00"
$wsStorage.Range("L11").NumberFormat = "0.00E+00"
$wsStorage.Range("F3:F9").NumberFormat = "0.0"

# Row 3
$wsStorage.Range("F3").Formula = "=20*100/1000"
$wsStorage.Range("J3").Value = 50
$wsStorage.Range("L3").Value = 10000000

# Row 4
$wsStorage.Range("F4").Formula = "=30*100/1000"
$wsStorage.Range("J4").Value = 60
$wsStorage.Range("L4").Value = 10000000

# Row 5
$wsStorage.Range("F5").Formula = "=40*100/1000"
$wsStorage.Range("J5").Value = 70
$wsStorage.Range("L5").Value = 10000000

# Row 6
$wsStorage.Range("F6").Formula = "=40*125/1000"
$wsStorage.Range("J6").Value = 80
$wsStorage.Range("L6").Value = 10000000

# Row 7
$wsStorage.Range("F7").Formula = "=40*150/1000"
$wsStorage.Range("J7").Value = 90
$wsStorage.Range("L7").Value = 10000000

# Row 8
$wsStorage.Range("F8").Formula = "=60*125/1000"
$wsStorage.Range("J8").Value = 100
$wsStorage.Range("L8").Value = 10000000

# Row 9
$wsStorage.Range("F9").Formula = "=60*150/1000"
$wsStorage.Range("J9").Value = 110
$wsStorage.Range("L9").Value = 10000000

# Row 11
$wsStorage.Range("J11").Value = 1000
$wsStorage.Range("L11").Value = 100000000

# ---------------------------------------------------------------------------
# Selections / active sheet: conversionAssets loses the tab focus, it moves
# to storageAssets. Set the non-active sheet's selection first so it isn't
# disturbed by the later activation.
# ---------------------------------------------------------------------------
$wsConversion.Range("F17").Select()
$wsStorage.Range("F8").Select()
